$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert the new benchmark row (row 7) mirroring the layout of rows 2-6:
#   A = instance name, D = MOVE TYPE, E = MAX_TRIALS label, F = MAX CANDIDATE, G = time limit label
$ws.Range("A7").Value = "fnl4461_n13380_uncorr_01.ttp"
$ws.Range("D7").Value = 5
$ws.Range("E7").Value = "num_cities/2"
$ws.Range("F7").Value = 5
$ws.Range("G7").Value = "DBL_MAX"

# Move the active selection to C8 (matches the post-edit cursor position)
$ws.Range("C8").Select() | Out-Null

# Configure the page for printing
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
